$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '68.311.26'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +2.99%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.641.55'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +2.36%  '

$ws.Range('E4').Value = '  -0.57%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '197.70'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +10.42%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '579.06'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.31%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.635.37'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +2.36%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.621'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.93%  '

$ws.Range('E9').Value = '  -0.49%  '

$ws.Range('E10').Value = '  +2.78%  '

$ws.Range('E11').Value = '  +10.13%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '56.75'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +7.34%  '

$ws.Range('E13').Value = '  +18.71%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.12'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.67%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.219.29'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.09%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.643.59'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.43%  '

$ws.Range('E17').Value = '  +0.73%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.60'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +5.19%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '68.260.17'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.90%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.65'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +3.38%  '

$ws.Range('E21').Value = '  +4.87%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '403.72'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +4.64%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '13.17'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +31.10%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.26'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.55%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '86.04'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.94%  '

$ws.Range('E26').Value = '  +4.97%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '12.66'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +5.50%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '3.88'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +8.85%  '

$ws.Range('E29').Value = '  +1.75%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.23'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +24.60%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '9.22'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +5.09%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '31.80'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.88%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '688.29'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +16.71%  '

$ws.Range('E34').Value = '  +4.80%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.118'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +7.01%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '64.64'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.41%  '

$ws.Range('E37').Value = '  +5.65%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.428'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +17.36%  '

$ws.Range('E39').Value = '  +0.13%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0₃0790'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +9.47%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.93'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +24.44%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.138'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +6.80%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.17'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +16.57%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.214.44'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +17.63%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.04'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +41.19%  '

$ws.Range('E46').Value = '  -0.32%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0422'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +4.99%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.99'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +10.82%  '

$ws.Range('E49').Value = '  +3.62%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.12'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.60%  '

$ws.Range('E51').Value = '  +3.91%  '
